# Automatische test-sync: 2025-07-23 22:49:50
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Logs" sheet: append a new row (row 26) with the 16th test mail entry.
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(26, 1).Value = "Kun jij dit even aan Koen doorgeven?"
$logs.Cells.Item(26, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(26, 3).Value = "Testmail #16: Kun jij dit even aan Koen doorgeven?"
$logs.Cells.Item(26, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item(26, 5).Value = "Geachte klant,`nBedankt voor uw e-mail. Kunt u ons wat meer context geven over het verzoek om iets aan Koen door te geven? Op deze manier kunnen we u beter van dienst zijn en ervoor zorgen dat uw bericht correct wordt doorgestuurd.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Cells.Item(26, 6).Value = "2025-07-23 22:49:28"
$logs.Cells.Item(26, 7).Value = "Ja"
$logs.Cells.Item(26, 8).Value = "Nee"
$logs.Cells.Item(26, 9).Value = "Ja"
$logs.Cells.Item(26, 10).Value = "Nee"

# Extend the conditional formatting ranges so row 26 is covered too
# (D, G, H, I, J columns), keeping all existing rules intact.
$dFc = $logs.Range("D2:D25").FormatConditions
for ($i = 1; $i -le $dFc.Count; $i++) {
    $dFc.Item($i).ModifyAppliesToRange($logs.Range("D2:D26"))
}

$gFc = $logs.Range("G2:G25").FormatConditions
for ($i = 1; $i -le $gFc.Count; $i++) {
    $gFc.Item($i).ModifyAppliesToRange($logs.Range("G2:G26"))
}

$hFc = $logs.Range("H2:H25").FormatConditions
for ($i = 1; $i -le $hFc.Count; $i++) {
    $hFc.Item($i).ModifyAppliesToRange($logs.Range("H2:H26"))
}

$iFc = $logs.Range("I2:I25").FormatConditions
for ($i = 1; $i -le $iFc.Count; $i++) {
    $iFc.Item($i).ModifyAppliesToRange($logs.Range("I2:I26"))
}

$jFc = $logs.Range("J2:J25").FormatConditions
for ($i = 1; $i -le $jFc.Count; $i++) {
    $jFc.Item($i).ModifyAppliesToRange($logs.Range("J2:J26"))
}

# ---------------------------------------------------------------------------
# 2. "Dashboard" sheet: swap rows 6/7 and append the new category row (11).
# ---------------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(6, 1).Value = "IT / Technisch probleem"
$dash.Cells.Item(7, 1).Value = "Factuur / Administratie"

$dash.Cells.Item(11, 1).Value = "Intern verzoek / Actie voor medewerker"
$dash.Cells.Item(11, 2).Value = 1

# ---------------------------------------------------------------------------
# 3. Chart on the Dashboard sheet: extend the category/value series ranges
#    from row 10 to row 11 to include the newly added category.
# ---------------------------------------------------------------------------
$chartObjects = $dash.ChartObjects()
$co = $chartObjects.Item(1)
$chart = $co.Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$11"
$series.Values = "='Dashboard'!`$B`$2:`$B`$11"
